$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)

$sh.TextFrame.TextRange.Text = "uBuy bietet den Kunden immer die besten personalisierten Angebote.`ruBuy gibt dem Benutzer ein außergewöhnliches Shopping Erlebnis."
